# Update leve-profit market data values across sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1117.2203
$ws.Range("I17").Value2 = 895.8889
$ws.Range("J17").Value2 = 1157.06
$ws.Range("K17").Value2 = 2687.6667
$ws.Range("L17").Value2 = 3471.18
$ws.Range("M17").Value2 = -2519.6667
$ws.Range("N17").Value2 = -3807.18

$ws.Range("H62").Value2 = 13232.81
$ws.Range("J62").Value2 = 5281.636
$ws.Range("L62").Value2 = 5281.636
$ws.Range("N62").Value2 = -6529.636

$ws.Range("H65").Value2 = 13232.81
$ws.Range("J65").Value2 = 5281.636
$ws.Range("L65").Value2 = 26408.18
$ws.Range("N65").Value2 = -32648.18

$ws.Range("H76").Value2 = 11467.1875
$ws.Range("J76").Value2 = 5291.1665
$ws.Range("L76").Value2 = 5291.1665
$ws.Range("N76").Value2 = -5921.1665

$ws.Range("H79").Value2 = 11467.1875
$ws.Range("J79").Value2 = 5291.1665
$ws.Range("L79").Value2 = 5291.1665
$ws.Range("N79").Value2 = -7475.1665

$ws.Range("H106").Value2 = 4078.4783
$ws.Range("I106").Value2 = 2857.5
$ws.Range("J106").Value2 = 5977.778
$ws.Range("K106").Value2 = 2857.5
$ws.Range("L106").Value2 = 5977.778
$ws.Range("M106").Value2 = -2226.5
$ws.Range("N106").Value2 = -7239.778

$ws.Range("H116").Value2 = 3627.7273
$ws.Range("I116").Value2 = 4443.294
$ws.Range("J116").Value2 = 2761.1875
$ws.Range("K116").Value2 = 4443.294
$ws.Range("L116").Value2 = 2761.1875
$ws.Range("M116").Value2 = -1001.294
$ws.Range("N116").Value2 = -9645.1875

$ws.Range("H128").Value2 = 35617.5
$ws.Range("J128").Value2 = 35617.5
$ws.Range("L128").Value2 = 35617.5
$ws.Range("N128").Value2 = -45577.5

$ws.Range("H129").Value2 = 777.8461
$ws.Range("I129").Value2 = 313.57144
$ws.Range("K129").Value2 = 940.71432
$ws.Range("M129").Value2 = 4059.28568

$ws.Range("H132").Value2 = 2019.7671
$ws.Range("I132").Value2 = 1119.7241
$ws.Range("J132").Value2 = 5499.933
$ws.Range("K132").Value2 = 3359.1723
$ws.Range("L132").Value2 = 16499.799
$ws.Range("M132").Value2 = -829.1722999999997
$ws.Range("N132").Value2 = -21559.799

$ws.Range("H137").Value2 = 2214.6545
$ws.Range("I137").Value2 = 2077.068
$ws.Range("K137").Value2 = 6231.204000000001
$ws.Range("M137").Value2 = -3681.204000000001

$ws.Range("H138").Value2 = 1726.6207
$ws.Range("I138").Value2 = 1218.5769
$ws.Range("J138").Value2 = 2139.4062
$ws.Range("K138").Value2 = 3655.7307
$ws.Range("L138").Value2 = 6418.2186
$ws.Range("M138").Value2 = 1484.2693
$ws.Range("N138").Value2 = -16698.2186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 370846.66
$ws.Range("I61").Value2 = 269888.72
$ws.Range("J61").Value2 = 629551.4
$ws.Range("K61").Value2 = 269888.72
$ws.Range("L61").Value2 = 629551.4
$ws.Range("M61").Value2 = -269676.72
$ws.Range("N61").Value2 = -629975.4

$ws.Range("H74").Value2 = 146274.44
$ws.Range("I74").Value2 = 173333.16
$ws.Range("J74").Value2 = 59085.223
$ws.Range("K74").Value2 = 173333.16
$ws.Range("L74").Value2 = 59085.223
$ws.Range("M74").Value2 = -172459.16
$ws.Range("N74").Value2 = -60833.223

$ws.Range("H77").Value2 = 146274.44
$ws.Range("I77").Value2 = 173333.16
$ws.Range("J77").Value2 = 59085.223
$ws.Range("K77").Value2 = 866665.8
$ws.Range("L77").Value2 = 295426.115
$ws.Range("M77").Value2 = -862297.8
$ws.Range("N77").Value2 = -304162.115

$ws.Range("H136").Value2 = 370846.66
$ws.Range("I136").Value2 = 269888.72
$ws.Range("J136").Value2 = 629551.4
$ws.Range("K136").Value2 = 809666.1599999999
$ws.Range("L136").Value2 = 1888654.2
$ws.Range("M136").Value2 = -807116.1599999999
$ws.Range("N136").Value2 = -1893754.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 639.51166
$ws.Range("I64").Value2 = 651.4167
$ws.Range("J64").Value2 = 624.4737
$ws.Range("K64").Value2 = 651.4167
$ws.Range("L64").Value2 = 624.4737
$ws.Range("M64").Value2 = -426.4167
$ws.Range("N64").Value2 = -1074.4737

$ws.Range("H67").Value2 = 639.51166
$ws.Range("I67").Value2 = 651.4167
$ws.Range("J67").Value2 = 624.4737
$ws.Range("K67").Value2 = 651.4167
$ws.Range("L67").Value2 = 624.4737
$ws.Range("M67").Value2 = 128.5833
$ws.Range("N67").Value2 = -2184.4737

$ws.Range("H109").Value2 = 21720.143
$ws.Range("J109").Value2 = 19736.666
$ws.Range("L109").Value2 = 19736.666
$ws.Range("N109").Value2 = -22510.666

$ws.Range("H134").Value2 = 2445.2769
$ws.Range("I134").Value2 = 2117.6382
$ws.Range("J134").Value2 = 3300.7778
$ws.Range("K134").Value2 = 6352.9146
$ws.Range("L134").Value2 = 9902.3334
$ws.Range("M134").Value2 = -3817.9146
$ws.Range("N134").Value2 = -14972.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 3601.2273
$ws.Range("I58").Value2 = 4172.2583
$ws.Range("J58").Value2 = 2239.5386
$ws.Range("K58").Value2 = 4172.2583
$ws.Range("L58").Value2 = 2239.5386
$ws.Range("M58").Value2 = -3969.2583
$ws.Range("N58").Value2 = -2645.5386

$ws.Range("H99").Value2 = 113323.445
$ws.Range("I99").Value2 = 126963.875
$ws.Range("J99").Value2 = 4200
$ws.Range("K99").Value2 = 126963.875
$ws.Range("L99").Value2 = 4200
$ws.Range("M99").Value2 = -125465.875
$ws.Range("N99").Value2 = -7196

$ws.Range("H126").Value2 = 113323.445
$ws.Range("I126").Value2 = 126963.875
$ws.Range("J126").Value2 = 4200
$ws.Range("K126").Value2 = 380891.625
$ws.Range("L126").Value2 = 12600
$ws.Range("M126").Value2 = -378421.625
$ws.Range("N126").Value2 = -17540

$ws.Range("H134").Value2 = 1479.7778
$ws.Range("I134").Value2 = 891.8919
$ws.Range("J134").Value2 = 2759.2942
$ws.Range("K134").Value2 = 2675.6757
$ws.Range("L134").Value2 = 8277.882599999999
$ws.Range("M134").Value2 = -140.6756999999998
$ws.Range("N134").Value2 = -13347.8826

$ws.Range("H136").Value2 = 3601.2273
$ws.Range("I136").Value2 = 4172.2583
$ws.Range("J136").Value2 = 2239.5386
$ws.Range("K136").Value2 = 12516.7749
$ws.Range("L136").Value2 = 6718.6158
$ws.Range("M136").Value2 = -9966.7749
$ws.Range("N136").Value2 = -11818.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 566.8570999999999
$ws.Range("I23").Value2 = 35
$ws.Range("J23").Value2 = 622.8421
$ws.Range("K23").Value2 = 105
$ws.Range("L23").Value2 = 1868.5263
$ws.Range("M23").Value2 = 130
$ws.Range("N23").Value2 = -2338.5263

$ws.Range("H68").Value2 = 358.18182
$ws.Range("I68").Value2 = 358
$ws.Range("J68").Value2 = 358.33334
$ws.Range("K68").Value2 = 1074
$ws.Range("L68").Value2 = 1075.00002
$ws.Range("M68").Value2 = -263
$ws.Range("N68").Value2 = -2697.00002

$ws.Range("H71").Value2 = 358.18182
$ws.Range("I71").Value2 = 358
$ws.Range("J71").Value2 = 358.33334
$ws.Range("K71").Value2 = 3222
$ws.Range("L71").Value2 = 3225.00006
$ws.Range("M71").Value2 = 834
$ws.Range("N71").Value2 = -11337.00006

$ws.Range("H86").Value2 = 596
$ws.Range("I86").Value2 = 596
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 1788
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = -602
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value2 = 596
$ws.Range("I89").Value2 = 596
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 5364
$ws.Range("L89").Value2 = 0
$ws.Range("M89").Value2 = 564
$ws.Range("N89").ClearContents()

$ws.Range("H127").Value2 = 1583.3334
$ws.Range("J127").Value2 = 1583.3334
$ws.Range("L127").Value2 = 4750.0002
$ws.Range("N127").Value2 = -14670.0002

$ws.Range("H131").Value2 = 1227.1774
$ws.Range("J131").Value2 = 1172.1698
$ws.Range("L131").Value2 = 3516.5094
$ws.Range("N131").Value2 = -13596.5094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 1099.9375
$ws.Range("I97").Value2 = 811
$ws.Range("J97").Value2 = 1581.5
$ws.Range("K97").Value2 = 811
$ws.Range("L97").Value2 = 1581.5
$ws.Range("M97").Value2 = -315
$ws.Range("N97").Value2 = -2573.5

$ws.Range("H132").Value2 = 2836.0193
$ws.Range("I132").Value2 = 2533.282
$ws.Range("J132").Value2 = 3744.2307
$ws.Range("K132").Value2 = 7599.846
$ws.Range("L132").Value2 = 11232.6921
$ws.Range("M132").Value2 = -5069.846
$ws.Range("N132").Value2 = -16292.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 4313
$ws.Range("I136").Value2 = 2413.2903
$ws.Range("J136").Value2 = 8239.066000000001
$ws.Range("K136").Value2 = 7239.8709
$ws.Range("L136").Value2 = 24717.198
$ws.Range("M136").Value2 = -4689.8709
$ws.Range("N136").Value2 = -29817.198

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value2 = 28000
$ws.Range("J63").Value2 = 28000
$ws.Range("L63").Value2 = 28000
$ws.Range("N63").Value2 = -29248

$ws.Range("H66").Value2 = 28000
$ws.Range("J66").Value2 = 28000
$ws.Range("L66").Value2 = 84000
$ws.Range("N66").Value2 = -90240

$ws.Range("H132").Value2 = 1433.0834
$ws.Range("I132").Value2 = 968.3333
$ws.Range("J132").Value2 = 2517.5
$ws.Range("K132").Value2 = 2904.9999
$ws.Range("L132").Value2 = 7552.5
$ws.Range("M132").Value2 = -374.9998999999998
$ws.Range("N132").Value2 = -12612.5

$ws.Range("H136").Value2 = 12483002
$ws.Range("I136").Value2 = 17562618
$ws.Range("J136").Value2 = 418915.03
$ws.Range("K136").Value2 = 52687854
$ws.Range("L136").Value2 = 1256745.09
$ws.Range("M136").Value2 = -52685304
$ws.Range("N136").Value2 = -1261845.09
